# Auto-generated edit script: applies cell-value updates produced by the
# scheduled Sheets runner, matching the target OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 749.4
$ws.Range("I18").Value = 749.4
$ws.Range("J18").Value = 0.0
$ws.Range("K18").Value = 749.4
$ws.Range("L18").Value = 0.0
$ws.Range("M18").Value = -465.4
$ws.Range("N18").ClearContents()
$ws.Range("H43").Value = 8063.625
$ws.Range("I43").Value = 2908.3333
$ws.Range("J43").Value = 11156.8
$ws.Range("K43").Value = 2908.3333
$ws.Range("L43").Value = 11156.8
$ws.Range("M43").Value = -2839.3333
$ws.Range("N43").Value = -11294.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4867.087
$ws.Range("I63").Value = 4951.9546
$ws.Range("K63").Value = 4951.9546
$ws.Range("M63").Value = -4265.9546
$ws.Range("H66").Value = 4867.087
$ws.Range("I66").Value = 4951.9546
$ws.Range("K66").Value = 24759.773
$ws.Range("M66").Value = -21327.773
$ws.Range("H81").Value = 0.0
$ws.Range("J81").Value = 0.0
$ws.Range("L81").Value = 0.0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0.0
$ws.Range("J84").Value = 0.0
$ws.Range("L84").Value = 0.0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 500.0
$ws.Range("J15").Value = 500.0
$ws.Range("L15").Value = 500.0
$ws.Range("N15").Value = -954.0
$ws.Range("H35").Value = 42325.0
$ws.Range("J35").Value = 42325.0
$ws.Range("L35").Value = 42325.0
$ws.Range("N35").Value = -42945.0
$ws.Range("H75").Value = 11523.777
$ws.Range("I75").Value = 6387.7144
$ws.Range("J75").Value = 29500.0
$ws.Range("K75").Value = 6387.7144
$ws.Range("L75").Value = 29500.0
$ws.Range("M75").Value = -5451.7144
$ws.Range("N75").Value = -31372.0
$ws.Range("H78").Value = 11523.777
$ws.Range("I78").Value = 6387.7144
$ws.Range("J78").Value = 29500.0
$ws.Range("K78").Value = 19163.1432
$ws.Range("L78").Value = 88500.0
$ws.Range("M78").Value = -14483.1432
$ws.Range("N78").Value = -97860.0
$ws.Range("H82").Value = 26059.092
$ws.Range("I82").Value = 11250.0
$ws.Range("J82").Value = 34521.43
$ws.Range("K82").Value = 11250.0
$ws.Range("L82").Value = 34521.43
$ws.Range("M82").Value = -10867.0
$ws.Range("N82").Value = -35287.43
$ws.Range("H85").Value = 26059.092
$ws.Range("I85").Value = 11250.0
$ws.Range("J85").Value = 34521.43
$ws.Range("K85").Value = 11250.0
$ws.Range("L85").Value = 34521.43
$ws.Range("M85").Value = -9924.0
$ws.Range("N85").Value = -37173.43

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 7000.0
$ws.Range("I41").Value = 7000.0
$ws.Range("K41").Value = 7000.0
$ws.Range("M41").Value = -6572.0
$ws.Range("H50").Value = 30000.0
$ws.Range("J50").Value = 30000.0
$ws.Range("L50").Value = 30000.0
$ws.Range("N50").Value = -31250.0
$ws.Range("H51").Value = 21316.666
$ws.Range("J51").Value = 21316.666
$ws.Range("L51").Value = 21316.666
$ws.Range("N51").Value = -22788.666
$ws.Range("H59").Value = 58679.0
$ws.Range("I59").Value = 45000.0
$ws.Range("J59").Value = 63238.668
$ws.Range("K59").Value = 45000.0
$ws.Range("L59").Value = 63238.668
$ws.Range("M59").Value = -43855.0
$ws.Range("N59").Value = -65528.668
$ws.Range("H60").Value = 15833.333
$ws.Range("I60").Value = 8000.0
$ws.Range("J60").Value = 19750.0
$ws.Range("K60").Value = 8000.0
$ws.Range("L60").Value = 19750.0
$ws.Range("M60").Value = -7489.0
$ws.Range("N60").Value = -20772.0
$ws.Range("H61").Value = 21316.666
$ws.Range("J61").Value = 21316.666
$ws.Range("L61").Value = 21316.666
$ws.Range("N61").Value = -22012.666
$ws.Range("H68").Value = 15598.333
$ws.Range("I68").Value = 5000.0
$ws.Range("J68").Value = 17718.0
$ws.Range("K68").Value = 5000.0
$ws.Range("L68").Value = 17718.0
$ws.Range("M68").Value = -4251.0
$ws.Range("N68").Value = -19216.0
$ws.Range("H71").Value = 15598.333
$ws.Range("I71").Value = 5000.0
$ws.Range("J71").Value = 17718.0
$ws.Range("K71").Value = 15000.0
$ws.Range("L71").Value = 53154.0
$ws.Range("M71").Value = -11256.0
$ws.Range("N71").Value = -60642.0
$ws.Range("H74").Value = 18979.264
$ws.Range("J74").Value = 18979.264
$ws.Range("L74").Value = 18979.264
$ws.Range("N74").Value = -20727.264
$ws.Range("H77").Value = 18979.264
$ws.Range("J77").Value = 18979.264
$ws.Range("L77").Value = 56937.792
$ws.Range("N77").Value = -65673.792
$ws.Range("H86").Value = 6617.909
$ws.Range("I86").Value = 4132.8335
$ws.Range("J86").Value = 9600.0
$ws.Range("K86").Value = 4132.8335
$ws.Range("L86").Value = 9600.0
$ws.Range("M86").Value = -3009.8335
$ws.Range("N86").Value = -11846.0
$ws.Range("H89").Value = 6617.909
$ws.Range("I89").Value = 4132.8335
$ws.Range("J89").Value = 9600.0
$ws.Range("K89").Value = 20664.1675
$ws.Range("L89").Value = 48000.0
$ws.Range("M89").Value = -15048.1675
$ws.Range("N89").Value = -59232.0
$ws.Range("H110").Value = 25999.666
$ws.Range("J110").Value = 25999.666
$ws.Range("L110").Value = 25999.666
$ws.Range("N110").Value = -34179.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 537.2759
$ws.Range("I5").Value = 332.61905
$ws.Range("K5").Value = 997.85715
$ws.Range("M5").Value = -885.85715
$ws.Range("H47").Value = 256.77777
$ws.Range("I47").Value = 251.57143
$ws.Range("J47").Value = 275.0
$ws.Range("K47").Value = 754.71429
$ws.Range("L47").Value = 825.0
$ws.Range("M47").Value = -323.71429
$ws.Range("N47").Value = -1687.0
$ws.Range("H75").Value = 1752.9412
$ws.Range("J75").Value = 1920.0
$ws.Range("L75").Value = 5760.0
$ws.Range("N75").Value = -7756.0
$ws.Range("H76").Value = 4986.0
$ws.Range("H78").Value = 1752.9412
$ws.Range("J78").Value = 1920.0
$ws.Range("L78").Value = 17280.0
$ws.Range("N78").Value = -27264.0
$ws.Range("H79").Value = 4986.0
$ws.Range("H94").Value = 5370.0
$ws.Range("J94").Value = 5370.0
$ws.Range("L94").Value = 16110.0
$ws.Range("N94").Value = -17462.0
$ws.Range("H112").Value = 4000.074
$ws.Range("I112").Value = 862.3333
$ws.Range("J112").Value = 4392.2915
$ws.Range("K112").Value = 2586.9999
$ws.Range("L112").Value = 13176.8745
$ws.Range("M112").Value = -1478.9999
$ws.Range("N112").Value = -15392.8745
$ws.Range("H131").Value = 879.3692
$ws.Range("I131").Value = 559.6667
$ws.Range("J131").Value = 911.88135
$ws.Range("K131").Value = 1679.0001
$ws.Range("L131").Value = 2735.64405
$ws.Range("M131").Value = 3360.9999
$ws.Range("N131").Value = -12815.64405
$ws.Range("H135").Value = 537.2759
$ws.Range("I135").Value = 332.61905
$ws.Range("K135").Value = 2993.57145
$ws.Range("M135").Value = -458.5714500000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4206.864
$ws.Range("I80").Value = 4461.6216
$ws.Range("J80").Value = 2860.2856
$ws.Range("K80").Value = 4461.6216
$ws.Range("L80").Value = 2860.2856
$ws.Range("M80").Value = -3463.6216
$ws.Range("N80").Value = -4856.2856
$ws.Range("H82").Value = 0.0
$ws.Range("J82").Value = 0.0
$ws.Range("L82").Value = 0.0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 4206.864
$ws.Range("I83").Value = 4461.6216
$ws.Range("J83").Value = 2860.2856
$ws.Range("K83").Value = 22308.108
$ws.Range("L83").Value = 14301.428
$ws.Range("M83").Value = -17316.108
$ws.Range("N83").Value = -24285.428
$ws.Range("H85").Value = 0.0
$ws.Range("J85").Value = 0.0
$ws.Range("L85").Value = 0.0
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 36189.0
$ws.Range("I87").Value = 0.0
$ws.Range("J87").Value = 36189.0
$ws.Range("K87").Value = 0.0
$ws.Range("L87").Value = 36189.0
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -38435.0
$ws.Range("H90").Value = 36189.0
$ws.Range("I90").Value = 0.0
$ws.Range("J90").Value = 36189.0
$ws.Range("K90").Value = 0.0
$ws.Range("L90").Value = 108567.0
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -119799.0

